# Generate Report for Handoff
# Adds a new row (row 9) for the file
# "df21252e-e450-4de9-8ca5-0a0f5034b9ea.md" to all three worksheets
# (Overview, zh-cn, de-de), mirroring the structure of the existing rows.

$wb = $excel.ActiveWorkbook

$fileId   = "df21252e-e450-4de9-8ca5-0a0f5034b9ea"
$zhToken  = "c602791f748c98086a4b73fdc7f232a94af041e9"

$mdName      = "$fileId.md"
$zhXlfName   = "$fileId.$zhToken.zh-cn.xlf"
$deXlfName   = "$fileId.$zhToken.de-de.xlf"

$overviewShaMd = "c3641b7552b5d2cf020eb2762f90bccba35cf923"
$zhShaMd       = "7a58e72762084bce882acbc6a519d12489903d60"
$zhShaXlf      = "2fc75185d36ff0a53dcfeb2d0008526d8af433c4"
$deShaMd       = "8bc2d839b3c5d2fa9a0a342338ae46acf6199695"
$deShaXlf      = "3781591e4e46dc70d9f234f5b8ecdb44ad8acf0f"

$overviewMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/$overviewShaMd/e2e/$mdName"
$zhMdUrl       = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/$zhShaMd/e2e/$mdName"
$zhXlfUrl      = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$zhShaXlf/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlfName"
$deMdUrl       = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/$deShaMd/e2e/$mdName"
$deXlfUrl      = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$deShaXlf/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlfName"

$readyForHandoff = "Ready for handoff"
$dotMd           = ".md"
$include         = "Include"
$epoch           = "0001-01-01 00:00:00"

$overviewDate = "2016-03-24 22:46:37"
$zhHandoffDate = "2016-03-24 22:46:30"

# Matches the existing "yyyy-mm-dd HH:mm:ss" number format (numFmtId 164)
# used by the other datetime cells (D2:D8, E2:E8, H2:H8, ...).
$dateTimeFormat = "yyyy-mm-dd HH:mm:ss"

# Matches the existing HyperLink font (underline, color FF6495ED) used by
# the other hyperlinked cells. Font.Color is BGR-encoded: 0x00ED9564.
$hyperlinkColor = 15570276

function Style-AsHyperlink($range) {
    $range.Font().Underline() = 2
    $range.Font().Color() = $hyperlinkColor
}

# ---------------------------------------------------------------------
# Sheet "Overview": new row 9
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B9").Value() = $readyForHandoff
$wsOverview.Range("C9").Value() = $readyForHandoff
$wsOverview.Range("D9").Value() = $overviewDate
$wsOverview.Range("D9").NumberFormat() = $dateTimeFormat

$wsOverview.Hyperlinks.Add($wsOverview.Range("A9"), $overviewMdUrl, "", "", $mdName)
Style-AsHyperlink($wsOverview.Range("A9"))

# ---------------------------------------------------------------------
# Sheet "zh-cn": new row 9
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B9").Value() = $dotMd
$wsZh.Range("C9").Value() = $readyForHandoff
$wsZh.Range("E9").Value() = $zhHandoffDate
$wsZh.Range("E9").NumberFormat() = $dateTimeFormat
$wsZh.Range("H9").Value() = $epoch
$wsZh.Range("H9").NumberFormat() = $dateTimeFormat
$wsZh.Range("J9").Value() = $include

$wsZh.Hyperlinks.Add($wsZh.Range("A9"), $zhMdUrl, "", "", $mdName)
Style-AsHyperlink($wsZh.Range("A9"))
$wsZh.Hyperlinks.Add($wsZh.Range("D9"), $zhXlfUrl, "", "", $zhXlfName)
Style-AsHyperlink($wsZh.Range("D9"))

# ---------------------------------------------------------------------
# Sheet "de-de": new row 9
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B9").Value() = $dotMd
$wsDe.Range("C9").Value() = $readyForHandoff
$wsDe.Range("E9").Value() = $overviewDate
$wsDe.Range("E9").NumberFormat() = $dateTimeFormat
$wsDe.Range("H9").Value() = $epoch
$wsDe.Range("H9").NumberFormat() = $dateTimeFormat
$wsDe.Range("J9").Value() = $include

$wsDe.Hyperlinks.Add($wsDe.Range("A9"), $deMdUrl, "", "", $mdName)
Style-AsHyperlink($wsDe.Range("A9"))
$wsDe.Hyperlinks.Add($wsDe.Range("D9"), $deXlfUrl, "", "", $deXlfName)
Style-AsHyperlink($wsDe.Range("D9"))
